# Applies the row permutation described in the commit diff to worksheet "Artfynd".
# Each data row (2-17) takes on the values that used to live in a different row of
# the same sheet (see the new_row <- old_row mapping in the comments below). Columns
# C, D, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AI, AT, AW, AX, AY are identical
# across every data row, so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- (content formerly in) row 11
$ws.Range("A2").Value = 111438447
$ws.Range("B2").Value = 76495
$ws.Range("E2").Value = 6487
$ws.Range("F2").Value = "Blågrå svartspik"
$ws.Range("G2").Value = "Chaenothecopsis fennica"
$ws.Range("H2").Value = "(Laurila) Tibell"
$ws.Range("Q2").Value = 468866.1318338988
$ws.Range("R2").Value = 6882808.390505624
$ws.Range("AO2").Value = "silverved tall"

# Row 3 <- (content formerly in) row 14
$ws.Range("A3").Value = 111438444
$ws.Range("B3").Value = 77515
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("Q3").Value = 468841.2361184616
$ws.Range("R3").Value = 6882806.276033297
$ws.Range("AO3").Value = "tall"

# Row 4 <- (content formerly in) row 8
$ws.Range("A4").Value = 111438439
$ws.Range("B4").Value = 78107
$ws.Range("E4").Value = 6453
$ws.Range("F4").Value = "Vedskivlav"
$ws.Range("G4").Value = "Hertelidea botryosa"
$ws.Range("H4").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q4").Value = 468788.4775288465
$ws.Range("R4").Value = 6882785.67140964
$ws.Range("AO4").Value = "silverved tall"

# Row 5 <- (content formerly in) row 17
$ws.Range("A5").Value = 111438446
$ws.Range("B5").Value = 77550
$ws.Range("E5").Value = 185
$ws.Range("F5").Value = "Violettgrå tagellav"
$ws.Range("G5").Value = "Bryoria nadvornikiana"
$ws.Range("H5").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("Q5").Value = 468853.3954624244
$ws.Range("R5").Value = 6882801.477506777
$ws.Range("AO5").Value = "gran"

# Row 6 <- (content formerly in) row 15
$ws.Range("A6").Value = 111438425
$ws.Range("B6").Value = 56398
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("Q6").Value = 468571.5178632676
$ws.Range("R6").Value = 6882722.999468728
$ws.Range("AO6").Value = "tall"
$ws.Range("M6").Value = "färska spår"
$ws.Range("AC6").Value = "Ringhack i tall"

# Row 7 <- (content formerly in) row 9
$ws.Range("A7").Value = 111438428
$ws.Range("B7").Value = 77597
$ws.Range("E7").Value = 864
$ws.Range("F7").Value = "Knottrig blåslav"
$ws.Range("G7").Value = "Hypogymnia bitteri"
$ws.Range("H7").Value = "(Lynge) Ahti"
$ws.Range("Q7").Value = 468740.5586073888
$ws.Range("R7").Value = 6882780.957796668
$ws.Range("AO7").Value = "björk"

# Row 8 <- (content formerly in) row 13
$ws.Range("A8").Value = 111438442
$ws.Range("B8").Value = 73696
$ws.Range("E8").Value = 6440
$ws.Range("F8").Value = "Vitgrynig nållav"
$ws.Range("G8").Value = "Chaenotheca subroscida"
$ws.Range("H8").Value = "(Eitner) Zahlbr."
$ws.Range("Q8").Value = 468800.3617588138
$ws.Range("R8").Value = 6882801.965499061
$ws.Range("AO8").Value = "gran"

# Row 9 <- (content formerly in) row 2
$ws.Range("A9").Value = 111438440
$ws.Range("B9").Value = 73689
$ws.Range("E9").Value = 308
$ws.Range("F9").Value = "Brunpudrad nållav"
$ws.Range("G9").Value = "Chaenotheca gracillima"
$ws.Range("H9").Value = "(Vain.) Tibell"
$ws.Range("Q9").Value = 468800.2970216064
$ws.Range("R9").Value = 6882794.936009536
$ws.Range("AO9").Value = "högstubbe björk"

# Row 10 <- (content formerly in) row 12
$ws.Range("A10").Value = 111438430
$ws.Range("B10").Value = 77268
$ws.Range("E10").Value = 228912
$ws.Range("F10").Value = "Mörk kolflarnlav"
$ws.Range("G10").Value = "Carbonicola myrmecina"
$ws.Range("H10").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q10").Value = 468756.5460031229
$ws.Range("R10").Value = 6882784.091042386
$ws.Range("AO10").Value = "brandstubbe"

# Row 11 <- (content formerly in) row 7
$ws.Range("A11").Value = 111438432
$ws.Range("B11").Value = 77267
$ws.Range("E11").Value = 6446
$ws.Range("F11").Value = "Kolflarnlav"
$ws.Range("G11").Value = "Carbonicola anthracophila"
$ws.Range("H11").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q11").Value = 468756.5460031229
$ws.Range("R11").Value = 6882784.091042386
$ws.Range("AO11").Value = "brandstubbe"

# Row 12 <- (content formerly in) row 4
$ws.Range("A12").Value = 111438457
$ws.Range("B12").Value = 78107
$ws.Range("E12").Value = 6453
$ws.Range("F12").Value = "Vedskivlav"
$ws.Range("G12").Value = "Hertelidea botryosa"
$ws.Range("H12").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q12").Value = 468747.5763832342
$ws.Range("R12").Value = 6882880.250689426
$ws.Range("AO12").Value = "silverved tall"

# Row 13 <- (content formerly in) row 5
$ws.Range("A13").Value = 111438435
$ws.Range("B13").Value = 76495
$ws.Range("E13").Value = 6487
$ws.Range("F13").Value = "Blågrå svartspik"
$ws.Range("G13").Value = "Chaenothecopsis fennica"
$ws.Range("H13").Value = "(Laurila) Tibell"
$ws.Range("Q13").Value = 468754.6686940129
$ws.Range("R13").Value = 6882784.108355919
$ws.Range("AO13").Value = "silverved tall"

# Row 14 <- (content formerly in) row 16
$ws.Range("A14").Value = 111438455
$ws.Range("B14").Value = 77515
$ws.Range("E14").Value = 6425
$ws.Range("F14").Value = "Garnlav"
$ws.Range("G14").Value = "Alectoria sarmentosa"
$ws.Range("H14").Value = "(Ach.) Ach."
$ws.Range("Q14").Value = 468784.2260541836
$ws.Range("R14").Value = 6882884.599394682
$ws.Range("AO14").Value = "tall"

# Row 15 <- (content formerly in) row 6
$ws.Range("A15").Value = 111438453
$ws.Range("B15").Value = 78107
$ws.Range("E15").Value = 6453
$ws.Range("F15").Value = "Vedskivlav"
$ws.Range("G15").Value = "Hertelidea botryosa"
$ws.Range("H15").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q15").Value = 468789.3971357156
$ws.Range("R15").Value = 6882885.489071017
$ws.Range("AO15").Value = "silverved tall"
$ws.Range("K15").ClearContents()
$ws.Range("L15").ClearContents()
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("AC15").ClearContents()

# Row 16 <- (content formerly in) row 3
$ws.Range("A16").Value = 111438433
$ws.Range("B16").Value = 78081
$ws.Range("E16").Value = 229821
$ws.Range("F16").Value = "Vedflamlav"
$ws.Range("G16").Value = "Ramboldia elabens"
$ws.Range("H16").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q16").Value = 468756.5460031229
$ws.Range("R16").Value = 6882784.091042386
$ws.Range("AO16").Value = "silverved tall"

# Row 17 <- (content formerly in) row 10
$ws.Range("A17").Value = 111438426
$ws.Range("B17").Value = 76918
$ws.Range("E17").Value = 6437
$ws.Range("F17").Value = "Blanksvart spiklav"
$ws.Range("G17").Value = "Calicium denigratum"
$ws.Range("H17").Value = "(Vain.) Tibell"
$ws.Range("Q17").Value = 468629.2461709682
$ws.Range("R17").Value = 6882722.464435354
$ws.Range("AO17").Value = "silverved tall"
